$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$abstractText = "Background`nid=`"Par1`">A COVID-19 outbreak started in Wuhan, China, last December and now has become a global pandemic.`n`n The clinical information in caring of critically ill patients with COVID-19 needs to be shared timely, especially under the situations that there is still a largely ongoing spread of COVID-19 in many countries.`n`n`nMethods`nid=`"Par2`">A multicenter prospective observational study investigated all the COVID-19 patients received in 19 ICUs of 16 hospitals in Wuhan, China, over 24 h between 8 AM February 2h and 8 AM February 27, 2020. The demographic information, clinical characteristics, vital signs, complications, laboratory values, and clinical managements of the patients were studied.`n`n`nResults`nid=`"Par3`">A total of 226 patients were included.`n`n Their median (interquartile range, IQR) age was 64 (57–70) years, and 139 (61.5%) patients were male.`n`n The duration from the date of ICU admission to the study date was 11 (5–17) days, and the duration from onset of symptoms to the study date was 31 (24–36) days.`n`n Among all the patients, 155 (68.6%) had at least one coexisting disease, and their sequential organ failure assessment score was 4 (2–8).`n`n Organ function damages were found in most of the patients: ARDS in 161 (71.2%) patients, septic shock in 34 (15.0%) patients, acute kidney injury occurred in 57 (25.2%) patients, cardiac injury in 61 (27.0%) patients, and lymphocytopenia in 160 (70.8%) patients.`n`n Of all the studied patients, 85 (37.6%) received invasive mechanical ventilation, including 14 (6.2%) treated with extracorporeal membrane oxygenation (ECMO) at the same time, 20 (8.8%) received noninvasive mechanical ventilation, and 24 (10.6%) received continuous renal replacement therapy.`n`n By April 9, 2020, 87 (38.5%) patients were deceased and 15 (6.7%) were still in the hospital.`n`n`nConclusions`nid=`"Par4`">Critically ill patients with COVID-19 are associated with a higher risk of severe complications and need to receive an intensive level of treatments.`n`n COVID-19 poses a great strain on critical care resources in hospitals.`n`n`nTrial registration`nid=`"Par5`">Chinese Clinical Trial Registry, ChiCTR2000030164. Registered on February 24, 2020, http://www.`n`nchictr.`n`norg.`n`ncn/edit.`n`naspx?pid=49983&amp;htm=4`n"
$authorsText = "[Yuan%Yu%NULL%0,   Dan%Xu%NULL%0,   Shouzhi%Fu%NULL%0,   Jun%Zhang%NULL%0,   Xiaobo%Yang%NULL%0,   Liang%Xu%NULL%0,   Jiqian%Xu%NULL%0,   Yongran%Wu%NULL%0,   Chaolin%Huang%NULL%0,   Yaqi%Ouyang%NULL%0,   Luyu%Yang%NULL%0,   Minghao%Fang%NULL%0,   Hongwen%Xiao%NULL%0,   Jing%Ma%NULL%0,   Wei%Zhu%NULL%0,   Song%Hu%NULL%0,   Quan%Hu%NULL%0,   Daoyin%Ding%NULL%0,   Ming%Hu%NULL%0,   Guochao%Zhu%NULL%0,   Weijiang%Xu%NULL%0,   Jun%Guo%NULL%0,   Jinglong%Xu%NULL%0,   Haitao%Yuan%NULL%0,   Bin%Zhang%NULL%0,   Zhui%Yu%yuzhui@whu.edu.cn%0,   Dechang%Chen%icudechangchen@163.com%0,   Shiying%Yuan%yuan_shiying@163.com%0,   You%Shang%you_shanghust@163.com%0]"

$ws.Range("D4").Value = $abstractText
$ws.Range("E4").Value = $authorsText
